$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "on Input relation"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "on Output relation"

# Header row
$ws2.Range("A1").Value = "Power Input"
$ws2.Range("B1").Value = "Hydrogen Output"
$ws2.Range("C1").Value = "Efficency "
$ws2.Range("E1").Value = "Piece 1"
$ws2.Range("F1").Value = "Piece 2"
$ws2.Range("G1").Value = "Piece 3"
$ws2.Range("H1").Value = "Total"
$ws2.Range("I1").Value = "Efficiency"
$ws2.Range("J1").Value = "Control"
$ws2.Range("K1").Value = "Identical"
$ws2.Range("M1").Value = "Range"
$ws2.Range("N1").Value = "Capacity Range"

# Data rows 2-25
# Row 2
$ws2.Range("A2").Value = 13.3
$ws2.Range("B2").Value = 10
$ws2.Range("C2").Formula = '=B2/A2'
$ws2.Range("E2").Formula = '=IF($B2>$N$2,$N$2*$P$2,B2*$O$2)'
$ws2.Range("F2").FormulaArray = '=IF(($B2-$N$2)>(0),_xlfn.IFS(B2>($N$2+$N$3),$N$3*$P$3,($B2-$N$2)>(0),(B2-$N$2)*$P$3),0)'
$ws2.Range("G2").Formula = '=IF(($B2-($N$2+$N$3))>(0),(B2-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H2").Formula = '=SUM(E2:G2)'
$ws2.Range("I2").Formula = '=H2/B2'
$ws2.Range("J2").Formula = '=1*$O$2'
$ws2.Range("K2").Formula = '=EXACT(J2,I2)'
$ws2.Range("M2").Value = 0.4
$ws2.Range("N2").Value = 15.037593984962406
$ws2.Range("O2").Value = 1.33
$ws2.Range("P2").Formula = '=O2'
# Row 3
$ws2.Range("A3").Value = 27.24812
$ws2.Range("B3").Value = 20
$ws2.Range("C3").Formula = '=B3/A3'
$ws2.Range("E3").Formula = '=IF($B3>$N$2,$N$2*$P$2,B3*$O$2)'
$ws2.Range("F3").FormulaArray = '=IF(($B3-$N$2)>(0),_xlfn.IFS(B3>($N$2+$N$3),$N$3*$P$3,($B3-$N$2)>(0),(B3-$N$2)*$P$3),0)'
$ws2.Range("G3").Formula = '=IF(($B3-($N$2+$N$3))>(0),(B3-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H3").Formula = '=SUM(E3:G3)'
$ws2.Range("I3").Formula = '=H3/B3'
$ws2.Range("J3").Formula = '=($N$2*$O$2+(B3-$N$2)*$O$3)/B3'
$ws2.Range("K3").Formula = '=EXACT(J3,I3)'
$ws2.Range("M3").Value = 0.7
$ws2.Range("N3").Value = 10.48951048951049
$ws2.Range("O3").Value = 1.44
$ws2.Range("P3").Formula = '=O3'
# Row 4
$ws2.Range("A4").Value = 34.479219999999998
$ws2.Range("B4").Value = 25
$ws2.Range("C4").Formula = '=B4/A4'
$ws2.Range("E4").Formula = '=IF($B4>$N$2,$N$2*$P$2,B4*$O$2)'
$ws2.Range("F4").FormulaArray = '=IF(($B4-$N$2)>(0),_xlfn.IFS(B4>($N$2+$N$3),$N$3*$P$3,($B4-$N$2)>(0),(B4-$N$2)*$P$3),0)'
$ws2.Range("G4").Formula = '=IF(($B4-($N$2+$N$3))>(0),(B4-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H4").Formula = '=SUM(E4:G4)'
$ws2.Range("I4").Formula = '=H4/B4'
$ws2.Range("J4").Formula = '=($N$2*$O$2+(B4-$N$2)*$O$3)/B4'
$ws2.Range("K4").Formula = '=EXACT(J4,I4)'
$ws2.Range("M4").Value = 1
$ws2.Range("N4").Value = 9.7402597402597397
$ws2.Range("O4").Value = 1.54
$ws2.Range("P4").Formula = '=O4'
# Row 5
$ws2.Range("A5").Value = 6.65
$ws2.Range("B5").Value = 5
$ws2.Range("C5").Formula = '=B5/A5'
$ws2.Range("E5").Formula = '=IF($B5>$N$2,$N$2*$P$2,B5*$O$2)'
$ws2.Range("F5").FormulaArray = '=IF(($B5-$N$2)>(0),_xlfn.IFS(B5>($N$2+$N$3),$N$3*$P$3,($B5-$N$2)>(0),(B5-$N$2)*$P$3),0)'
$ws2.Range("G5").Formula = '=IF(($B5-($N$2+$N$3))>(0),(B5-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H5").Formula = '=SUM(E5:G5)'
$ws2.Range("I5").Formula = '=H5/B5'
$ws2.Range("J5").Formula = '=1*$O$2'
$ws2.Range("K5").Formula = '=EXACT(J5,I5)'
# Row 6
$ws2.Range("A6").Value = 13.3
$ws2.Range("B6").Value = 10
$ws2.Range("C6").Formula = '=B6/A6'
$ws2.Range("E6").Formula = '=IF($B6>$N$2,$N$2*$P$2,B6*$O$2)'
$ws2.Range("F6").FormulaArray = '=IF(($B6-$N$2)>(0),_xlfn.IFS(B6>($N$2+$N$3),$N$3*$P$3,($B6-$N$2)>(0),(B6-$N$2)*$P$3),0)'
$ws2.Range("G6").Formula = '=IF(($B6-($N$2+$N$3))>(0),(B6-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H6").Formula = '=SUM(E6:G6)'
$ws2.Range("I6").Formula = '=H6/B6'
$ws2.Range("J6").Formula = '=1*$O$2'
$ws2.Range("K6").Formula = '=EXACT(J6,I6)'
# Row 7
$ws2.Range("A7").Value = 20.048120000000001
$ws2.Range("B7").Value = 15
$ws2.Range("C7").Formula = '=B7/A7'
$ws2.Range("E7").Formula = '=IF($B7>$N$2,$N$2*$P$2,B7*$O$2)'
$ws2.Range("F7").FormulaArray = '=IF(($B7-$N$2)>(0),_xlfn.IFS(B7>($N$2+$N$3),$N$3*$P$3,($B7-$N$2)>(0),(B7-$N$2)*$P$3),0)'
$ws2.Range("G7").Formula = '=IF(($B7-($N$2+$N$3))>(0),(B7-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H7").Formula = '=SUM(E7:G7)'
$ws2.Range("I7").Formula = '=H7/B7'
$ws2.Range("J7").Formula = '=1*$O$2'
$ws2.Range("K7").Formula = '=EXACT(J7,I7)'
# Row 8
$ws2.Range("A8").Value = 27.24812
$ws2.Range("B8").Value = 20
$ws2.Range("C8").Formula = '=B8/A8'
$ws2.Range("E8").Formula = '=IF($B8>$N$2,$N$2*$P$2,B8*$O$2)'
$ws2.Range("F8").FormulaArray = '=IF(($B8-$N$2)>(0),_xlfn.IFS(B8>($N$2+$N$3),$N$3*$P$3,($B8-$N$2)>(0),(B8-$N$2)*$P$3),0)'
$ws2.Range("G8").Formula = '=IF(($B8-($N$2+$N$3))>(0),(B8-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H8").Formula = '=SUM(E8:G8)'
$ws2.Range("I8").Formula = '=H8/B8'
$ws2.Range("J8").Formula = '=($N$2*$O$2+(B8-$N$2)*$O$3)/B8'
$ws2.Range("K8").Formula = '=EXACT(J8,I8)'
# Row 9
$ws2.Range("A9").Value = 42.179220000000001
$ws2.Range("B9").Value = 30
$ws2.Range("C9").Formula = '=B9/A9'
$ws2.Range("E9").Formula = '=IF($B9>$N$2,$N$2*$P$2,B9*$O$2)'
$ws2.Range("F9").FormulaArray = '=IF(($B9-$N$2)>(0),_xlfn.IFS(B9>($N$2+$N$3),$N$3*$P$3,($B9-$N$2)>(0),(B9-$N$2)*$P$3),0)'
$ws2.Range("G9").Formula = '=IF(($B9-($N$2+$N$3))>(0),(B9-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H9").Formula = '=SUM(E9:G9)'
$ws2.Range("I9").Formula = '=H9/B9'
$ws2.Range("J9").Formula = '=($N$2*$O$2+($N$3)*$O$3+(B9-$N$2-$N$3)*$O$4)/B9'
$ws2.Range("K9").Formula = '=EXACT(J9,I9)'
# Row 10
$ws2.Range("A10").Value = 49.879219999999997
$ws2.Range("B10").Value = 35
$ws2.Range("C10").Formula = '=B10/A10'
$ws2.Range("E10").Formula = '=IF($B10>$N$2,$N$2*$P$2,B10*$O$2)'
$ws2.Range("F10").FormulaArray = '=IF(($B10-$N$2)>(0),_xlfn.IFS(B10>($N$2+$N$3),$N$3*$P$3,($B10-$N$2)>(0),(B10-$N$2)*$P$3),0)'
$ws2.Range("G10").Formula = '=IF(($B10-($N$2+$N$3))>(0),(B10-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H10").Formula = '=SUM(E10:G10)'
$ws2.Range("I10").Formula = '=H10/B10'
$ws2.Range("J10").Formula = '=($N$2*$O$2+($N$3)*$O$3+(B10-$N$2-$N$3)*$O$4)/B10'
$ws2.Range("K10").Formula = '=EXACT(J10,I10)'
# Row 11
$ws2.Range("A11").Value = 45.259219999999999
$ws2.Range("B11").Value = 32
$ws2.Range("C11").Formula = '=B11/A11'
$ws2.Range("E11").Formula = '=IF($B11>$N$2,$N$2*$P$2,B11*$O$2)'
$ws2.Range("F11").FormulaArray = '=IF(($B11-$N$2)>(0),_xlfn.IFS(B11>($N$2+$N$3),$N$3*$P$3,($B11-$N$2)>(0),(B11-$N$2)*$P$3),0)'
$ws2.Range("G11").Formula = '=IF(($B11-($N$2+$N$3))>(0),(B11-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H11").Formula = '=SUM(E11:G11)'
$ws2.Range("I11").Formula = '=H11/B11'
$ws2.Range("J11").Formula = '=($N$2*$O$2+($N$3)*$O$3+(B11-$N$2-$N$3)*$O$4)/B11'
$ws2.Range("K11").Formula = '=EXACT(J11,I11)'
# Row 12
$ws2.Range("A12").Value = 36.019219999999997
$ws2.Range("B12").Value = 26
$ws2.Range("C12").Formula = '=B12/A12'
$ws2.Range("E12").Formula = '=IF($B12>$N$2,$N$2*$P$2,B12*$O$2)'
$ws2.Range("F12").FormulaArray = '=IF(($B12-$N$2)>(0),_xlfn.IFS(B12>($N$2+$N$3),$N$3*$P$3,($B12-$N$2)>(0),(B12-$N$2)*$P$3),0)'
$ws2.Range("G12").Formula = '=IF(($B12-($N$2+$N$3))>(0),(B12-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H12").Formula = '=SUM(E12:G12)'
$ws2.Range("I12").Formula = '=H12/B12'
$ws2.Range("J12").Formula = '=($N$2*$O$2+($N$3)*$O$3+(B12-$N$2-$N$3)*$O$4)/B12'
$ws2.Range("K12").Formula = '=EXACT(J12,I12)'
# Row 13
$ws2.Range("A13").Value = 33.008119999999998
$ws2.Range("B13").Value = 24
$ws2.Range("C13").Formula = '=B13/A13'
$ws2.Range("E13").Formula = '=IF($B13>$N$2,$N$2*$P$2,B13*$O$2)'
$ws2.Range("F13").FormulaArray = '=IF(($B13-$N$2)>(0),_xlfn.IFS(B13>($N$2+$N$3),$N$3*$P$3,($B13-$N$2)>(0),(B13-$N$2)*$P$3),0)'
$ws2.Range("G13").Formula = '=IF(($B13-($N$2+$N$3))>(0),(B13-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H13").Formula = '=SUM(E13:G13)'
$ws2.Range("I13").Formula = '=H13/B13'
$ws2.Range("J13").Formula = '=($N$2*$O$2+(B13-$N$2)*$O$3)/B13'
$ws2.Range("K13").Formula = '=EXACT(J13,I13)'
# Row 14
$ws2.Range("A14").Value = 24.368120000000001
$ws2.Range("B14").Value = 18
$ws2.Range("C14").Formula = '=B14/A14'
$ws2.Range("E14").Formula = '=IF($B14>$N$2,$N$2*$P$2,B14*$O$2)'
$ws2.Range("F14").FormulaArray = '=IF(($B14-$N$2)>(0),_xlfn.IFS(B14>($N$2+$N$3),$N$3*$P$3,($B14-$N$2)>(0),(B14-$N$2)*$P$3),0)'
$ws2.Range("G14").Formula = '=IF(($B14-($N$2+$N$3))>(0),(B14-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H14").Formula = '=SUM(E14:G14)'
$ws2.Range("I14").Formula = '=H14/B14'
$ws2.Range("J14").Formula = '=($N$2*$O$2+(B14-$N$2)*$O$3)/B14'
$ws2.Range("K14").Formula = '=EXACT(J14,I14)'
# Row 15
$ws2.Range("A15").Value = 5.32
$ws2.Range("B15").Value = 4
$ws2.Range("C15").Formula = '=B15/A15'
$ws2.Range("E15").Formula = '=IF($B15>$N$2,$N$2*$P$2,B15*$O$2)'
$ws2.Range("F15").FormulaArray = '=IF(($B15-$N$2)>(0),_xlfn.IFS(B15>($N$2+$N$3),$N$3*$P$3,($B15-$N$2)>(0),(B15-$N$2)*$P$3),0)'
$ws2.Range("G15").Formula = '=IF(($B15-($N$2+$N$3))>(0),(B15-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H15").Formula = '=SUM(E15:G15)'
$ws2.Range("I15").Formula = '=H15/B15'
$ws2.Range("J15").Formula = '=1*$O$2'
$ws2.Range("K15").Formula = '=EXACT(J15,I15)'
# Row 16
$ws2.Range("A16").Value = 13.3
$ws2.Range("B16").Value = 10
$ws2.Range("C16").Formula = '=B16/A16'
$ws2.Range("E16").Formula = '=IF($B16>$N$2,$N$2*$P$2,B16*$O$2)'
$ws2.Range("F16").FormulaArray = '=IF(($B16-$N$2)>(0),_xlfn.IFS(B16>($N$2+$N$3),$N$3*$P$3,($B16-$N$2)>(0),(B16-$N$2)*$P$3),0)'
$ws2.Range("G16").Formula = '=IF(($B16-($N$2+$N$3))>(0),(B16-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H16").Formula = '=SUM(E16:G16)'
$ws2.Range("I16").Formula = '=H16/B16'
$ws2.Range("J16").Formula = '=1*$O$2'
$ws2.Range("K16").Formula = '=EXACT(J16,I16)'
# Row 17
$ws2.Range("A17").Value = 34.479219999999998
$ws2.Range("B17").Value = 25
$ws2.Range("C17").Formula = '=B17/A17'
$ws2.Range("E17").Formula = '=IF($B17>$N$2,$N$2*$P$2,B17*$O$2)'
$ws2.Range("F17").FormulaArray = '=IF(($B17-$N$2)>(0),_xlfn.IFS(B17>($N$2+$N$3),$N$3*$P$3,($B17-$N$2)>(0),(B17-$N$2)*$P$3),0)'
$ws2.Range("G17").Formula = '=IF(($B17-($N$2+$N$3))>(0),(B17-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H17").Formula = '=SUM(E17:G17)'
$ws2.Range("I17").Formula = '=H17/B17'
$ws2.Range("J17").Formula = '=($N$2*$O$2+(B17-$N$2)*$O$3)/B17'
$ws2.Range("K17").Formula = '=EXACT(J17,I17)'
# Row 18
$ws2.Range("A18").Value = 34.479219999999998
$ws2.Range("B18").Value = 25
$ws2.Range("C18").Formula = '=B18/A18'
$ws2.Range("E18").Formula = '=IF($B18>$N$2,$N$2*$P$2,B18*$O$2)'
$ws2.Range("F18").FormulaArray = '=IF(($B18-$N$2)>(0),_xlfn.IFS(B18>($N$2+$N$3),$N$3*$P$3,($B18-$N$2)>(0),(B18-$N$2)*$P$3),0)'
$ws2.Range("G18").Formula = '=IF(($B18-($N$2+$N$3))>(0),(B18-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H18").Formula = '=SUM(E18:G18)'
$ws2.Range("I18").Formula = '=H18/B18'
$ws2.Range("J18").Formula = '=($N$2*$O$2+(B18-$N$2)*$O$3)/B18'
$ws2.Range("K18").Formula = '=EXACT(J18,I18)'
# Row 19
$ws2.Range("A19").Value = 42.179220000000001
$ws2.Range("B19").Value = 30
$ws2.Range("C19").Formula = '=B19/A19'
$ws2.Range("E19").Formula = '=IF($B19>$N$2,$N$2*$P$2,B19*$O$2)'
$ws2.Range("F19").FormulaArray = '=IF(($B19-$N$2)>(0),_xlfn.IFS(B19>($N$2+$N$3),$N$3*$P$3,($B19-$N$2)>(0),(B19-$N$2)*$P$3),0)'
$ws2.Range("G19").Formula = '=IF(($B19-($N$2+$N$3))>(0),(B19-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H19").Formula = '=SUM(E19:G19)'
$ws2.Range("I19").Formula = '=H19/B19'
$ws2.Range("J19").Formula = '=($N$2*$O$2+($N$3)*$O$3+(B19-$N$2-$N$3)*$O$4)/B19'
$ws2.Range("K19").Formula = '=EXACT(J19,I19)'
# Row 20
$ws2.Range("A20").Value = 27.24812
$ws2.Range("B20").Value = 20
$ws2.Range("C20").Formula = '=B20/A20'
$ws2.Range("E20").Formula = '=IF($B20>$N$2,$N$2*$P$2,B20*$O$2)'
$ws2.Range("F20").FormulaArray = '=IF(($B20-$N$2)>(0),_xlfn.IFS(B20>($N$2+$N$3),$N$3*$P$3,($B20-$N$2)>(0),(B20-$N$2)*$P$3),0)'
$ws2.Range("G20").Formula = '=IF(($B20-($N$2+$N$3))>(0),(B20-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H20").Formula = '=SUM(E20:G20)'
$ws2.Range("I20").Formula = '=H20/B20'
$ws2.Range("J20").Formula = '=($N$2*$O$2+(B20-$N$2)*$O$3)/B20'
$ws2.Range("K20").Formula = '=EXACT(J20,I20)'
# Row 21
$ws2.Range("A21").Value = 13.3
$ws2.Range("B21").Value = 10
$ws2.Range("C21").Formula = '=B21/A21'
$ws2.Range("E21").Formula = '=IF($B21>$N$2,$N$2*$P$2,B21*$O$2)'
$ws2.Range("F21").FormulaArray = '=IF(($B21-$N$2)>(0),_xlfn.IFS(B21>($N$2+$N$3),$N$3*$P$3,($B21-$N$2)>(0),(B21-$N$2)*$P$3),0)'
$ws2.Range("G21").Formula = '=IF(($B21-($N$2+$N$3))>(0),(B21-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H21").Formula = '=SUM(E21:G21)'
$ws2.Range("I21").Formula = '=H21/B21'
$ws2.Range("J21").Formula = '=1*$O$2'
$ws2.Range("K21").Formula = '=EXACT(J21,I21)'
# Row 22
$ws2.Range("A22").Value = 6.65
$ws2.Range("B22").Value = 5
$ws2.Range("C22").Formula = '=B22/A22'
$ws2.Range("E22").Formula = '=IF($B22>$N$2,$N$2*$P$2,B22*$O$2)'
$ws2.Range("F22").FormulaArray = '=IF(($B22-$N$2)>(0),_xlfn.IFS(B22>($N$2+$N$3),$N$3*$P$3,($B22-$N$2)>(0),(B22-$N$2)*$P$3),0)'
$ws2.Range("G22").Formula = '=IF(($B22-($N$2+$N$3))>(0),(B22-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H22").Formula = '=SUM(E22:G22)'
$ws2.Range("I22").Formula = '=H22/B22'
$ws2.Range("J22").Formula = '=1*$O$2'
$ws2.Range("K22").Formula = '=EXACT(J22,I22)'
# Row 23
$ws2.Range("A23").Value = 42.179220000000001
$ws2.Range("B23").Value = 30
$ws2.Range("C23").Formula = '=B23/A23'
$ws2.Range("E23").Formula = '=IF($B23>$N$2,$N$2*$P$2,B23*$O$2)'
$ws2.Range("F23").FormulaArray = '=IF(($B23-$N$2)>(0),_xlfn.IFS(B23>($N$2+$N$3),$N$3*$P$3,($B23-$N$2)>(0),(B23-$N$2)*$P$3),0)'
$ws2.Range("G23").Formula = '=IF(($B23-($N$2+$N$3))>(0),(B23-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H23").Formula = '=SUM(E23:G23)'
$ws2.Range("I23").Formula = '=H23/B23'
$ws2.Range("J23").Formula = '=($N$2*$O$2+($N$3)*$O$3+(B23-$N$2-$N$3)*$O$4)/B23'
$ws2.Range("K23").Formula = '=EXACT(J23,I23)'
# Row 24
$ws2.Range("A24").Value = 34.479219999999998
$ws2.Range("B24").Value = 25
$ws2.Range("C24").Formula = '=B24/A24'
$ws2.Range("E24").Formula = '=IF($B24>$N$2,$N$2*$P$2,B24*$O$2)'
$ws2.Range("F24").FormulaArray = '=IF(($B24-$N$2)>(0),_xlfn.IFS(B24>($N$2+$N$3),$N$3*$P$3,($B24-$N$2)>(0),(B24-$N$2)*$P$3),0)'
$ws2.Range("G24").Formula = '=IF(($B24-($N$2+$N$3))>(0),(B24-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H24").Formula = '=SUM(E24:G24)'
$ws2.Range("I24").Formula = '=H24/B24'
$ws2.Range("J24").Formula = '=($N$2*$O$2+(B24-$N$2)*$O$3)/B24'
$ws2.Range("K24").Formula = '=EXACT(J24,I24)'
# Row 25
$ws2.Range("A25").Value = 24.368120000000001
$ws2.Range("B25").Value = 18
$ws2.Range("C25").Formula = '=B25/A25'
$ws2.Range("E25").Formula = '=IF($B25>$N$2,$N$2*$P$2,B25*$O$2)'
$ws2.Range("F25").FormulaArray = '=IF(($B25-$N$2)>(0),_xlfn.IFS(B25>($N$2+$N$3),$N$3*$P$3,($B25-$N$2)>(0),(B25-$N$2)*$P$3),0)'
$ws2.Range("G25").Formula = '=IF(($B25-($N$2+$N$3))>(0),(B25-($N$2+$N$3))*$P$4,0)'
$ws2.Range("H25").Formula = '=SUM(E25:G25)'
$ws2.Range("I25").Formula = '=H25/B25'
$ws2.Range("J25").Formula = '=($N$2*$O$2+(B25-$N$2)*$O$3)/B25'
$ws2.Range("K25").Formula = '=EXACT(J25,I25)'

# Selection to match target sheet view
[void]$ws2.Range("N35").Select()